$wb = $excel.ActiveWorkbook

# --- pitstop sheet: add tyre_before / tyre_after columns ---
$ws = $wb.Worksheets.Item("pitstop")
$ws.Range("H1").Value = "tyre_before"
$ws.Range("I1").Value = "tyre_after"

$tyreData = @(
    @(4,3),
    @(5,3),
    @(5,3),
    @(4,3),
    @(4,3),
    @(3,5),
    @(3,4),
    @(4,3),
    @(5,3),
    @(4,3),
    @(3,5),
    @(5,5),
    @(4,3),
    @(3,5),
    @(5,3),
    @(4,3),
    @(3,5),
    @(5,3),
    @(3,4),
    @(5,3),
    @(5,3),
    @(5,3),
    @(3,4),
    @(4,3),
    @(5,3),
    @(4,3),
    @(3,4),
    @(4,3),
    @(3,5)
)

for ($i = 0; $i -lt $tyreData.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $tyreData[$i][0]
    $ws.Cells.Item($row, 9).Value = $tyreData[$i][1]
}

$ws.Activate()
$ws.Range("H31").Select()

# --- weather sheet (race-day weather conditions) ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$weather = $wb.Worksheets.Add($null, $lastSheet)
$weather.Name = "weather"

$weather.Range("A1").Value = "Skycondition"
$weather.Range("B1").Value = "Humid and Mostly Cloudy"
$weather.Range("A2").Value = "Temperature"
$weather.Range("B2").Value = "85.56°F"
$weather.Range("A3").Value = "Humidity"
$weather.Range("B3").Value = 0.7
$weather.Range("B3").NumberFormat = "0%"
$weather.Range("A4").Value = "Wind speed"
$weather.Range("B4").Value = "11.61 mph"
$weather.Range("A5").Value = "Wind bearing"
$weather.Range("B5").Value = "139°"

$weather.Activate()
$weather.Range("G11").Select()

# --- altitude sheet (circuit altitude delta) ---
$altitude = $wb.Worksheets.Add($null, $weather)
$altitude.Name = "altitude"
$altitude.Range("A1").Value = "delta"
$altitude.Range("B1").Value = 5.3

$altitude.Activate()
$altitude.Range("B2").Select()

Write-Host "Done"
